$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (B5:AH5, columns 2..34): round the existing measured values to 2
# decimal places (custom accuracy).
for ($col = 2; $col -le 34; $col++) {
    $cell = $ws.Cells.Item(5, $col)
    $cur = $cell.Value2
    if ($cur -ne $null) {
        $cell.Value = [Math]::Round([double]$cur, 2)
    }
}

# Drop the extra data row (row 6) entirely - shrinks the used range back
# down to A1:AH5.
$ws.Rows.Item(6).Delete()
